$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> "Data effettiva avanzamento" (column Z) date serial value
$zValues = @{
    2 = 45926
    3 = 45926
    4 = 45926
    5 = 45926
    6 = 45924
    7 = 45923
    8 = 45922
    9 = 45922
    10 = 45922
    11 = 45922
    12 = 45919
    13 = 45919
    14 = 45918
    15 = 45917
    16 = 45916
    17 = 45916
    18 = 45915
    19 = 45933
    20 = 45915
    21 = 45915
    22 = 45915
    23 = 45915
    24 = 45912
    25 = 45911
    26 = 45911
    27 = 45909
    28 = 45909
    29 = 45909
    30 = 45909
    31 = 45909
    32 = 45909
    33 = 45908
    34 = 45908
    35 = 45908
    36 = 45937
    37 = 45898
    38 = 45897
    39 = 45897
    40 = 45897
    41 = 45897
    42 = 45897
    43 = 45863
    44 = 45896
    45 = 45887
    46 = 45887
    47 = 45887
    48 = 45887
    49 = 45887
    50 = 45887
    51 = 45876
    52 = 45874
    53 = 45873
    54 = 45870
    55 = 45870
    56 = 45869
    57 = 45869
    58 = 45869
    59 = 45869
    60 = 45867
    61 = 45866
    62 = 45866
    63 = 45866
    64 = 45866
    65 = 45863
    66 = 45862
    67 = 45860
    68 = 45860
    69 = 45860
    70 = 45859
    71 = 45859
    72 = 45855
    73 = 45854
    74 = 45854
    75 = 45854
    76 = 45853
    77 = 45853
    78 = 45852
    79 = 45852
    80 = 45852
    81 = 45852
    82 = 45852
    83 = 45852
    84 = 45849
    85 = 45848
    86 = 45845
    87 = 45841
    88 = 45840
    89 = 45839
    90 = 45839
    91 = 45839
}

foreach ($row in $zValues.Keys) {
    $cell = $ws.Cells.Item($row, 26)  # Column Z = 26
    $cell.Value = $zValues[$row]
    $cell.NumberFormat = "YYYY-MM-DD"
}

